$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/19/2024  Through  8/25/2024"

# --- Type conversions (copy style/type from stable template cells in row 14, then set value) ---
# numeric -> text "0" (shared string index 20), style 14
$ws.Range("D14").Copy($ws.Range("C15"))
$ws.Range("D14").Copy($ws.Range("C27"))
$ws.Range("D14").Copy($ws.Range("F29"))
$ws.Range("D14").Copy($ws.Range("F30"))
$ws.Range("D14").Copy($ws.Range("D33"))

# numeric -> text "***.*" (shared string index 21), style 14
$ws.Range("E14").Copy($ws.Range("E33"))

# text -> numeric, style 15
$ws.Range("I14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("I14").Copy($ws.Range("C33"))
$ws.Range("C33").Value = 1
$ws.Range("I14").Copy($ws.Range("F33"))
$ws.Range("F33").Value = 1

# text -> numeric, style 16
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100

# --- Plain value updates ---
$ws.Range("N14").Value = -85.714285714285
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = -22.222222222222
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -21.052631578947
$ws.Range("I16").Value = 89
$ws.Range("J16").Value = 111
$ws.Range("K16").Value = -19.819819819819
$ws.Range("L16").Value = -27.049180327868
$ws.Range("M16").Value = -58.604651162790
$ws.Range("N16").Value = -86.696562032884
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -36.111111111111
$ws.Range("I17").Value = 283
$ws.Range("J17").Value = 298
$ws.Range("K17").Value = -5.033557046979
$ws.Range("L17").Value = -3.741496598639
$ws.Range("M17").Value = 46.632124352331
$ws.Range("N17").Value = 8.429118773946
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -31.578947368421
$ws.Range("I18").Value = 121
$ws.Range("J18").Value = 143
$ws.Range("K18").Value = -15.384615384615
$ws.Range("L18").Value = 5.217391304347
$ws.Range("M18").Value = -48.290598290598
$ws.Range("N18").Value = -88.113948919449
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = -63.636363636363
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 85
$ws.Range("H19").Value = -49.411764705882
$ws.Range("I19").Value = 380
$ws.Range("J19").Value = 460
$ws.Range("K19").Value = -17.391304347826
$ws.Range("L19").Value = -10.377358490566
$ws.Range("M19").Value = 32.404181184669
$ws.Range("N19").Value = -2.061855670103
$ws.Range("C20").Value = 13
$ws.Range("E20").Value = 8.333333333333
$ws.Range("F20").Value = 44
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = 22.222222222222
$ws.Range("I20").Value = 263
$ws.Range("J20").Value = 217
$ws.Range("K20").Value = 21.198156682027
$ws.Range("L20").Value = 68.589743589743
$ws.Range("M20").Value = 6.477732793522
$ws.Range("N20").Value = -88.040018190086
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 51
$ws.Range("E21").Value = -35.294117647058
$ws.Range("F21").Value = 141
$ws.Range("G21").Value = 196
$ws.Range("H21").Value = -28.061224489795
$ws.Range("I21").Value = 1159
$ws.Range("J21").Value = 1243
$ws.Range("K21").Value = -6.757843925985
$ws.Range("L21").Value = 3.297682709447
$ws.Range("M21").Value = -3.336113427856
$ws.Range("N21").Value = -74.672202797202
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -45.454545454545
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 110
$ws.Range("H24").Value = -26.363636363636
$ws.Range("I24").Value = 737
$ws.Range("J24").Value = 817
$ws.Range("K24").Value = -9.791921664626
$ws.Range("L24").Value = -26.593625498008
$ws.Range("M24").Value = 35.229357798165
$ws.Range("C25").Value = 8
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 9.090909090909
$ws.Range("I25").Value = 177
$ws.Range("J25").Value = 138
$ws.Range("K25").Value = 28.260869565217
$ws.Range("L25").Value = -25.316455696202
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 15.384615384615
$ws.Range("F26").Value = 62
$ws.Range("G26").Value = 54
$ws.Range("H26").Value = 14.814814814814
$ws.Range("I26").Value = 518
$ws.Range("J26").Value = 417
$ws.Range("K26").Value = 24.220623501199
$ws.Range("L26").Value = 21.028037383177
$ws.Range("M26").Value = 24.519230769230
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 150
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = 39.130434782608
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 37
$ws.Range("K28").Value = 23.333333333333
$ws.Range("L28").Value = 0
$ws.Range("M29").Value = -63.333333333333
$ws.Range("N29").Value = -73.809523809523
$ws.Range("M30").Value = -52.380952380952
$ws.Range("N30").Value = -75.609756097561
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = -50
$ws.Range("I33").Value = 2
$ws.Range("K33").Value = -71.428571428571
$ws.Range("L33").Value = -60
